$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update goal values for experiment A
$ws.Range("B12").Value = 275
$ws.Range("B13").Value = -550
$ws.Range("B17").Value = -0.551
$ws.Range("B18").Value = 0.835

# Update the selected cell to match the new active selection
$ws.Range("B18").Select()
